$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for dates that were dropped from the report (holiday/weekend
# placeholder rows with near-zero activity): 12/25 (r57), 12/31 (r63), 1/1 (r64),
# 1/3 (r66), 1/4 (r67). Deleting bottom-to-top keeps the remaining row indices valid.
$ws.Rows("67:67").Delete()
$ws.Rows("66:66").Delete()
$ws.Rows("64:64").Delete()
$ws.Rows("63:63").Delete()
$ws.Rows("57:57").Delete()

# Refresh the remaining rows (now 57-63) with the corrected figures from the
# re-uploaded source data.
$data = @(
    @(46017, 5611, 3085, 2829, 176, 45, 31, 3, 1),
    @(46018, 5611, 25, 25, 0, 0, 0, 0, 0),
    @(46019, 5611, 29, 28, 1, 0, 0, 0, 0),
    @(46020, 5609, 3495, 3223, 203, 39, 25, 4, 1),
    @(46021, 5606, 3525, 3264, 197, 40, 20, 4, 0),
    @(46024, 5596, 3159, 2895, 191, 38, 32, 2, 1),
    @(46027, 5587, 4160, 3861, 216, 48, 31, 4, 0)
)

$arr = New-Object 'object[,]' 7,9
for ($r = 0; $r -lt 7; $r++) {
    for ($c = 0; $c -lt 9; $c++) {
        $arr[$r,$c] = $data[$r][$c]
    }
}
$ws.Range("A57:I63").Value = $arr

# The trailing rows (64-68) stay blank, but column A keeps the date number
# format that was carried down the column.
$ws.Range("A64:A68").NumberFormat = "d-mmm-yy"

# Match the author's final selection (row 63, the new last data row).
[void]$ws.Range("A63:I63").Select()
